$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 6 new blank rows at 905 (existing rows 905-945 shift down to 911-951)
$ws.Range("A905:R910").EntireRow.Insert()

# Copy the formatting (incl. date number format on column D) from the row
# just below the inserted block (old row 905, now at 911) onto the new rows.
$ws.Range("A911:R911").Copy()
$ws.Range("A905:R910").PasteSpecial(-4122)  # xlPasteFormats

$const = @{
    A = 1
    B = 'Agrícola del Norte S.A. de Arica'
    C = 'Arica y Parinacota'
    E = 15
    F = 100112020
    G = 'Tomate'
    H = 'Larga vida'
    O = 'Región de Arica y Parinacota'
    R = 'Hortaliza'
}

$newRows = @(
    @{ Row=905; I='Primera'; J=250; K=6500; L=7000; M=6750; N='$/bandeja 18 kilos'; P=375; Q=18 },
    @{ Row=906; I='Primera'; J=350; K=3000; L=3500; M=3250; N='$/caja 10 kilos';    P=325; Q=10 },
    @{ Row=907; I='Segunda'; J=375; K=6000; L=6500; M=6200; N='$/bandeja 18 kilos'; P=344; Q=18 },
    @{ Row=908; I='Segunda'; J=450; K=2500; L=3000; M=2750; N='$/caja 10 kilos';    P=275; Q=10 },
    @{ Row=909; I='Tercera'; J=450; K=5000; L=5500; M=5250; N='$/bandeja 18 kilos'; P=292; Q=18 },
    @{ Row=910; I='Tercera'; J=450; K=2000; L=2500; M=2250; N='$/caja 10 kilos';    P=225; Q=10 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Range("A$row").Value = $const.A
    $ws.Range("B$row").Value = $const.B
    $ws.Range("C$row").Value = $const.C
    $ws.Range("D$row").Value = Get-Date -Year 2022 -Month 7 -Day 11 -Hour 0 -Minute 0 -Second 0
    $ws.Range("E$row").Value = $const.E
    $ws.Range("F$row").Value = $const.F
    $ws.Range("G$row").Value = $const.G
    $ws.Range("H$row").Value = $const.H
    $ws.Range("I$row").Value = $r.I
    $ws.Range("J$row").Value = $r.J
    $ws.Range("K$row").Value = $r.K
    $ws.Range("L$row").Value = $r.L
    $ws.Range("M$row").Value = $r.M
    $ws.Range("N$row").Value = $r.N
    $ws.Range("O$row").Value = $const.O
    $ws.Range("P$row").Value = $r.P
    $ws.Range("Q$row").Value = $r.Q
    $ws.Range("R$row").Value = $const.R
}
